$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update "last updated" timestamp (row 1) ---
$ws.Range("A1").Value = "Datos actualizados a 23 de Agosto de 2020 a las 20:47"

# --- Country rank swaps (text labels) ---
# Siria overtook Guinea-Bisau in total cases: rows 135/136 swap countries
$ws.Range("A135").Value = "Siria"
$ws.Range("A136").Value = "Guinea-Bisau"

# Timor Oriental overtook Santa Lucia: rows 202/203 swap countries
$ws.Range("A202").Value = "Timor Oriental"
$ws.Range("A203").Value = "Santa Lucia"

# --- Updated statistics per country row ---
# Row 4
$ws.Range("B4").Value = 5856961
$ws.Range("C4").Value = 15533
$ws.Range("D4").Value = 3153012
$ws.Range("E4").Value = 2523559
$ws.Range("G4").Value = 216
$ws.Range("H4").Value = 180390

# Row 6
$ws.Range("B6").Value = 3105038
$ws.Range("C6").Value = 61602
$ws.Range("D6").Value = 2336763
$ws.Range("E6").Value = 710584
$ws.Range("G6").Value = 845
$ws.Range("H6").Value = 57691

# Row 21
$ws.Range("B21").Value = 258249
$ws.Range("C21").Value = 1217
$ws.Range("D21").Value = 237165
$ws.Range("E21").Value = 14963
$ws.Range("G21").Value = 19
$ws.Range("H21").Value = 6121

# Row 22
$ws.Range("B22").Value = 242899
$ws.Range("C22").Value = 4897
$ws.Range("E22").Value = 127436
$ws.Range("G22").Value = 1
$ws.Range("H22").Value = 30513

# Row 23
$ws.Range("B23").Value = 234290
$ws.Range("C23").Value = 433
$ws.Range("E23").Value = 16008

# Row 27
$ws.Range("B27").Value = 124893
$ws.Range("C27").Value = 264
$ws.Range("D27").Value = 111098
$ws.Range("E27").Value = 4722
$ws.Range("G27").Value = 2
$ws.Range("H27").Value = 9073

# Row 30
$ws.Range("B30").Value = 107769
$ws.Range("C30").Value = 680
$ws.Range("D30").Value = 94816
$ws.Range("E30").Value = 6643
$ws.Range("G30").Value = 33
$ws.Range("H30").Value = 6310

# Row 58
$ws.Range("B58").Value = 41460
$ws.Range("C58").Value = 392
$ws.Range("D58").Value = 29142
$ws.Range("E58").Value = 10883
$ws.Range("G58").Value = 11
$ws.Range("H58").Value = 1435

# Row 74
$ws.Range("B74").Value = 21867
$ws.Range("C74").Value = 77
$ws.Range("D74").Value = 16119
$ws.Range("E74").Value = 5336
$ws.Range("G74").Value = 1
$ws.Range("H74").Value = 412

# Row 75
$ws.Range("D75").Value = 11843
$ws.Range("E75").Value = 6831
$ws.Range("G75").Value = 3
$ws.Range("H75").Value = 128

# Row 135
$ws.Range("B135").Value = 2217
$ws.Range("C135").Value = 74
$ws.Range("D135").Value = 505
$ws.Range("E135").Value = 1623
$ws.Range("G135").Value = 4
$ws.Range("H135").Value = 89

# Row 136
$ws.Range("B136").Value = 2149
$ws.Range("D136").Value = 1015
$ws.Range("E136").Value = 1101
$ws.Range("H136").Value = 33

# Row 141
$ws.Range("B141").Value = 1911
$ws.Range("C141").Value = 4
$ws.Range("D141").Value = 1086
$ws.Range("E141").Value = 272
$ws.Range("G141").Value = 7
$ws.Range("H141").Value = 553

# Row 161
$ws.Range("B161").Value = 986
$ws.Range("C161").Value = 4
$ws.Range("D161").Value = 870
$ws.Range("E161").Value = 40
